{"js": "// Replace the answer text in every cell of the single table, left-to-right /\n// top-to-bottom (row-major order), with the new set of arithmetic problems.\n// We overwrite each cell's first paragraph via insertText(..., replace) so the\n// existing run/paragraph formatting (TimeNewRoman, sz 30, left-justified) is\n// preserved -- only the <w:t> content changes, matching the source diff.\nconst NEW_VALUES = [\"18+22=40\", \"83-81=2\", \"2+42=44\", \"18+12=30\", \"87+2=89\", \"51-15=36\", \"28-23=5\", \"66-32=34\", \"77-2=75\", \"81-45=36\", \"6+53=59\", \"69+10=79\", \"97-42=55\", \"59+14=73\", \"40+55=95\", \"80-28=52\", \"78-35=43\", \"44-28=16\", \"19+26=45\", \"22-4=18\", \"12+24=36\", \"52-6=46\", \"69-57=12\", \"70-14=56\", \"11+65=76\", \"83-9=74\", \"67-37=30\", \"86-3=83\", \"64-31=33\", \"55-38=17\", \"27+21=48\", \"92-10=82\", \"64+3=67\", \"97-5=92\", \"49+4=53\", \"27+15=42\", \"5+59=64\", \"28+34=62\", \"25-9=16\", \"77-28=49\", \"51+16=67\", \"85-41=44\", \"98-96=2\", \"72-41=31\", \"90-72=18\", \"27+35=62\", \"56-46=10\", \"54-41=13\", \"98-76=22\", \"54-2=52\", \"61+22=83\", \"79-38=41\", \"65-46=19\", \"9+0=9\", \"79+18=97\", \"36+46=82\", \"68-64=4\", \"25+41=66\", \"58-36=22\", \"20+3=23\", \"94-14=80\", \"71-43=28\", \"61-2=59\", \"15+37=52\", \"26+25=51\", \"37-16=21\", \"53+34=87\", \"26+39=65\", \"17+24=41\", \"71-36=35\", \"33+19=52\", \"5+28=33\", \"39-28=11\", \"47+20=67\", \"32+43=75\", \"97-9=88\", \"37+15=52\", \"87-36=51\", \"51+19=70\", \"14+72=86\", \"22+63=85\", \"51-16=35\", \"81-79=2\", \"74+7=81\", \"55+36=91\", \"84-0=84\", \"36+4=40\", \"54+25=79\", \"6+31=37\", \"52-50=2\", \"63-1=62\", \"56-7=49\", \"56-35=21\", \"29+62=91\", \"53-13=40\", \"71-34=37\", \"16+25=41\", \"15-0=15\", \"23-20=3\", \"70-4=66\"];\n\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst colCount = NEW_VALUES.length / rowCount;\n\n// Collect the first paragraph of every cell, in row-major order.\nconst paragraphs = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const cell = table.getCell(r, c);\n    paragraphs.push(cell.body.paragraphs.getFirst());\n  }\n}\n\n// Replace each paragraph's text while preserving its run formatting.\nfor (let i = 0; i < paragraphs.length; i++) {\n  paragraphs[i].insertText(NEW_VALUES[i], Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the answer text in every cell of the single table, left-to-right /\n# top-to-bottom (row-major order), with the new set of arithmetic problems.\n# Only the visible text changes -- cell/paragraph/run formatting is left as-is\n# because we write straight into Cell.Range.Text instead of re-creating the\n# paragraph.\n$newValues = @(\n    \"18+22=40\",\n    \"83-81=2\",\n    \"2+42=44\",\n    \"18+12=30\",\n    \"87+2=89\",\n    \"51-15=36\",\n    \"28-23=5\",\n    \"66-32=34\",\n    \"77-2=75\",\n    \"81-45=36\",\n    \"6+53=59\",\n    \"69+10=79\",\n    \"97-42=55\",\n    \"59+14=73\",\n    \"40+55=95\",\n    \"80-28=52\",\n    \"78-35=43\",\n    \"44-28=16\",\n    \"19+26=45\",\n    \"22-4=18\",\n    \"12+24=36\",\n    \"52-6=46\",\n    \"69-57=12\",\n    \"70-14=56\",\n    \"11+65=76\",\n    \"83-9=74\",\n    \"67-37=30\",\n    \"86-3=83\",\n    \"64-31=33\",\n    \"55-38=17\",\n    \"27+21=48\",\n    \"92-10=82\",\n    \"64+3=67\",\n    \"97-5=92\",\n    \"49+4=53\",\n    \"27+15=42\",\n    \"5+59=64\",\n    \"28+34=62\",\n    \"25-9=16\",\n    \"77-28=49\",\n    \"51+16=67\",\n    \"85-41=44\",\n    \"98-96=2\",\n    \"72-41=31\",\n    \"90-72=18\",\n    \"27+35=62\",\n    \"56-46=10\",\n    \"54-41=13\",\n    \"98-76=22\",\n    \"54-2=52\",\n    \"61+22=83\",\n    \"79-38=41\",\n    \"65-46=19\",\n    \"9+0=9\",\n    \"79+18=97\",\n    \"36+46=82\",\n    \"68-64=4\",\n    \"25+41=66\",\n    \"58-36=22\",\n    \"20+3=23\",\n    \"94-14=80\",\n    \"71-43=28\",\n    \"61-2=59\",\n    \"15+37=52\",\n    \"26+25=51\",\n    \"37-16=21\",\n    \"53+34=87\",\n    \"26+39=65\",\n    \"17+24=41\",\n    \"71-36=35\",\n    \"33+19=52\",\n    \"5+28=33\",\n    \"39-28=11\",\n    \"47+20=67\",\n    \"32+43=75\",\n    \"97-9=88\",\n    \"37+15=52\",\n    \"87-36=51\",\n    \"51+19=70\",\n    \"14+72=86\",\n    \"22+63=85\",\n    \"51-16=35\",\n    \"81-79=2\",\n    \"74+7=81\",\n    \"55+36=91\",\n    \"84-0=84\",\n    \"36+4=40\",\n    \"54+25=79\",\n    \"6+31=37\",\n    \"52-50=2\",\n    \"63-1=62\",\n    \"56-7=49\",\n    \"56-35=21\",\n    \"29+62=91\",\n    \"53-13=40\",\n    \"71-34=37\",\n    \"16+25=41\",\n    \"15-0=15\",\n    \"23-20=3\",\n    \"70-4=66\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$idx]\n    $idx++\n  }\n}\n"}
